$d = $word.ActiveDocument

# "NTT Address" block on the first page used to be split across two
# paragraphs:
#   Cluj Napoca
#   City, State 400158
# Merge them into a single paragraph reading "Cluj Napoca City, 400158"
# (the trailing paragraph mark/paragraph is removed along with it).
$paraMark = [char]13
$searchText = "Cluj Napoca" + $paraMark + "City, State 400158"

$found = $d.Content.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "Cluj Napoca City, 400158", 2)

if (-not $found) {
    throw "Could not find the NTT address paragraphs to merge"
}
